# v 0.5 - added board caching to improve move calculation speed
# Adds a "Time 8" results column (K) with the new cached-move-generation
# timings for the first several opening moves, and documents the change
# in the notes list below the summary tables.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New column header
$ws.Range("K2").Value = "Time 8"

# New timing values for the Time 8 run (only the first few moves were
# re-measured)
$ws.Range("K3").Value  = 0.59
$ws.Range("K4").Value  = 0.47
$ws.Range("K5").Value  = 0.56
$ws.Range("K6").Value  = 1.18
$ws.Range("K7").Value  = 0.86
$ws.Range("K8").Value  = 0.55
$ws.Range("K9").Value  = 1.52
$ws.Range("K10").Value = 0.64
$ws.Range("K11").Value = 0.72

# Document the change alongside the other "Time N = ..." notes
$ws.Range("A49").Value = "Time 8 = added transposition table / move cache"

# Put the active selection on the newly added column (matches the author's
# saved cursor position)
[void]$ws.Range("K12").Select()
